# "change in serial squence" — extend the serial-packet layout diagrams on
# the sheet: the top frame (row 1/2) gains a "minestate" + "termination"
# byte, and the bottom frame (row 7/8) is rebuilt with an extra byte and
# the "right encoder"/"z angle" fields swapped/resized.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# Row 1 — byte index strip 1..12 (content unchanged, only re-striped
# colours change, handled in the styling pass below)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = 1
$ws.Range("B1").Formula = "=A1+1"
$ws.Range("C1:L1").Formula = "=B1+1"

# ---------------------------------------------------------------------
# Row 2 — field labels. K2 becomes "minestate" (was "termination"),
# "termination" moves out to the newly added L2.
# ---------------------------------------------------------------------
$ws.Range("K2").Value = "minestate"
$ws.Range("L2").Value = "termination"

# ---------------------------------------------------------------------
# Row 7 — byte index strip, now 1..13 and carrying real formulas like
# row 1 (it used to be bare literal values).
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("B7").Formula = "=A7+1"
$ws.Range("C7:J7").Formula = "=B7+1"
$ws.Range("K7").Formula = "=J7+1"
$ws.Range("L7:M7").Formula = "=K7+1"

# ---------------------------------------------------------------------
# Row 8 — field labels, rebuilt: "left encoder" widens to A8:D8,
# "right encoder" moves to E8:H8, "z angle" moves to I8:K8,
# "minestate" moves from J8 to L8, "termination" newly added at M8.
# ---------------------------------------------------------------------
$ws.Range("D8:F8").UnMerge()
$ws.Range("G8:I8").UnMerge()
$ws.Range("A8:C8").UnMerge()

$ws.Range("D8").Value = $null
$ws.Range("G8").Value = $null
$ws.Range("J8").Value = $null

$ws.Range("A8:D8").Merge()
$ws.Range("E8").Value = "right encoder"
$ws.Range("E8:H8").Merge()
$ws.Range("I8").Value = "z angle"
$ws.Range("I8:K8").Merge()
$ws.Range("L8").Value = "minestate"
$ws.Range("M8").Value = "termination"

# ---------------------------------------------------------------------
# Row 15 — a few blank formatted cells left under the widened columns.
# ---------------------------------------------------------------------
$ws.Range("L15:N15").WrapText = $false

# ---------------------------------------------------------------------
# New column widths for the two newly-meaningful columns L:M.
# ---------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 14.43
$ws.Columns.Item(13).ColumnWidth = 14.43

# ---------------------------------------------------------------------
# Styling pass — restripe Good/Bad/Normal cell styles + alignment to
# match the new layout. Style is (re)applied before alignment, since
# assigning a named Style resets direct alignment formatting.
# ---------------------------------------------------------------------

function Set-GoodCentered($rng) {
    $r = $ws.Range($rng)
    $r.Style = "Good"
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
}
function Set-BadCentered($rng) {
    $r = $ws.Range($rng)
    $r.Style = "Bad"
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
}
function Set-NormalCentered($rng) {
    $r = $ws.Range($rng)
    $r.Style = "Normal"
    $r.HorizontalAlignment = $xlCenter
}

# Row 1
Set-GoodCentered("A1:B1")
Set-BadCentered("C1")
Set-BadCentered("D1:F1")
Set-BadCentered("G1:J1")
Set-BadCentered("K1")
Set-BadCentered("L1")

# Row 2
Set-GoodCentered("A2:C2")
Set-BadCentered("D2:F2")
Set-GoodCentered("G2:J2")
Set-BadCentered("K2")
Set-GoodCentered("L2")

# Row 5 (unchanged visually, keep Normal + horizontal-center)
Set-NormalCentered("D5:F5")

# Row 7
Set-GoodCentered("A7:B7")
Set-BadCentered("C7:D7")
Set-GoodCentered("E7:H7")
Set-BadCentered("I7:K7")
Set-GoodCentered("L7")
Set-BadCentered("M7")

# Row 8
Set-GoodCentered("A8:D8")
Set-BadCentered("E8:H8")
Set-GoodCentered("I8:K8")
Set-GoodCentered("L8")
Set-GoodCentered("M8")

# ---------------------------------------------------------------------
# Dimension / selection housekeeping (Excel normally derives the used
# range automatically, but make the final selection match explicitly).
# ---------------------------------------------------------------------
$ws.Range("G14").Select()
